# dlc13_iec61400-1ed3.xlsx — update wind-direction offset from +/-10 deg to +/-8 deg
# (wdir350 -> wdir352, i.e. F = -10 -> -8; wdir010 -> wdir008, i.e. F = 10 -> 8).
# The "Case id." column (D) is a shared formula that derives its text from
# E/F/G, so it recalculates automatically once F is updated - no need to
# touch D directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$lastRow = $ws.Range("F1").End(-4121).Row   # xlDown

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 6)   # column F
    $v = $cell.Value2
    if ($v -eq -10) {
        $cell.Value = -8
    } elseif ($v -eq 10) {
        $cell.Value = 8
    }
}
